$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F6").Value = 100
    $ws.Range("F11").Value = 36
    $ws.Range("F12").Value = 318
    $ws.Range("F16").Value = 102
    $ws.Range("F17").Value = 16
    $ws.Range("F20").Value = 102
    $ws.Range("F21").Value = 1003
    $ws.Range("F22").Value = 1417
    $ws.Range("F23").Value = 309
    $ws.Range("F26").Value = 81
    $ws.Range("F30").Value = 232
    $ws.Range("F31").Value = 262
    $ws.Range("F32").Value = 282
    $ws.Range("F33").Value = 1641
    $ws.Range("F34").Value = 55
    $ws.Range("F37").Value = 591
    $ws.Range("F39").Value = 3776
    $ws.Range("F40").Value = 440
    $ws.Range("F41").Value = 213
    $ws.Range("F42").Value = 931
    $ws.Range("F45").Value = 78
    $ws.Range("F46").Value = 43
}
